$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: shift row1 data from B1:M1 to A1:L1 (preserves styles incl. empty styled cells)
$ws.Range("B1:M1").Copy($ws.Range("A1"))
$ws.Range("H1").ClearContents()
$ws.Range("M1").Clear()

# Step 2: fix G1 value (was shifted-in 21213, should become 1)
$ws.Range("G1").Value = 1

# Step 3: build row 2 as a copy of row1's A:G block, then L
$ws.Range("A1:G1").Copy($ws.Range("A2"))
$ws.Range("L1").Copy($ws.Range("L2"))
$ws.Range("G2").Value = 2

# Step 4: selection / view state
[void]$ws.Range("B1:L2").Select()
